# "Informacja o graczu i wprowadzilem zmienna sleep"
#
# Slide 3 ("Co w tym czasie sie dzialo?") carries the running stats about
# the TicTacToe project.  Adding player info + a sleep variable grew the
# source file, so the cached line-count bumps from 268 to 281 lines, and
# the screenshot illustrating the code was nudged to the right to make
# room for the extra text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# --- "Koncowa liczba linijek kodu: 268" -> "...: 281" -------------------
$statsBox = $s.Shapes.Item(3)
$statsRange = $statsBox.TextFrame.TextRange
$lineCountParagraph = $statsRange.Paragraphs(3)
$lineCountRun = $lineCountParagraph.Runs(2)
$lineCountRun.Text = "281"

# --- Shift the accompanying screenshot to the right ----------------------
$picture = $s.Shapes.Item(4)
$emuPerPoint = 914400 / 72
$newLeftEmu = 8269174
$picture.Left = $newLeftEmu / $emuPerPoint
